# Apply updated ANT/OVA PSSM & PSERM scores (supplemental figures/tables recalculation)
# D = ANT PSSM Score, E = ANT PSERM Score, H = OVA PSSM Score, I = OVA PSERM Score
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 18.70835811617782
$ws.Range("E2").Value = 2.394884232768993
$ws.Range("H2").Value = 16.61188720856606
$ws.Range("I2").Value = 0.2984133251572303
$ws.Range("D3").Value = 21.78781819507843
$ws.Range("E3").Value = 3.794802146959263
$ws.Range("H3").Value = 18.56074126420972
$ws.Range("I3").Value = 0.5677252160905606
$ws.Range("D4").Value = 21.21466642386854
$ws.Range("E4").Value = 4.080868320297949
$ws.Range("H4").Value = 17.33500891643433
$ws.Range("I4").Value = 0.2012108128637404
$ws.Range("D5").Value = 20.93574573440618
$ws.Range("E5").Value = 3.210347172685364
$ws.Range("H5").Value = 18.22564475672517
$ws.Range("I5").Value = 0.5002461950043608
$ws.Range("D6").Value = 14.03055252409274
$ws.Range("E6").Value = -1.401431320986156
$ws.Range("H6").Value = 15.62993622491712
$ws.Range("I6").Value = 0.1979523798382264
$ws.Range("D7").Value = 20.25757979524642
$ws.Range("E7").Value = 2.87724787262922
$ws.Range("H7").Value = 17.89301979498221
$ws.Range("I7").Value = 0.512687872365001
$ws.Range("D8").Value = 16.28660367971138
$ws.Range("E8").Value = 1.65304917474928
$ws.Range("H8").Value = 15.28955857419547
$ws.Range("I8").Value = 0.6560040692333706
$ws.Range("D9").Value = 10.13251572283791
$ws.Range("E9").Value = -5.503575924106495
$ws.Range("H9").Value = 15.07107894917308
$ws.Range("I9").Value = -0.5650126977713252
$ws.Range("D10").Value = 14.6118160766635
$ws.Range("E10").Value = -0.7511575072227321
$ws.Range("H10").Value = 15.90878608759149
$ws.Range("I10").Value = 0.5458125037052646
$ws.Range("D11").Value = 22.77725557143994
$ws.Range("E11").Value = 4.494211496700024
$ws.Range("H11").Value = 18.55104612863017
$ws.Range("I11").Value = 0.2680020538902441
$ws.Range("D12").Value = 22.81044396872054
$ws.Range("E12").Value = 4.546661955161667
$ws.Range("H12").Value = 18.74742774369592
$ws.Range("I12").Value = 0.4836457301370458
$ws.Range("D13").Value = 20.45184264837387
$ws.Range("E13").Value = 3.471031157598676
$ws.Range("H13").Value = 17.42187231657475
$ws.Range("I13").Value = 0.4410608257995616
$ws.Range("D14").Value = 12.22179047109426
$ws.Range("E14").Value = -4.272327153196745
$ws.Range("H14").Value = 15.4042377241739
$ws.Range("I14").Value = -1.089879900117107
$ws.Range("D15").Value = 21.15629313688745
$ws.Range("E15").Value = 2.833080809257746
$ws.Range("H15").Value = 18.68700455831132
$ws.Range("I15").Value = 0.3637922306816095
$ws.Range("D16").Value = 18.46045835614044
$ws.Range("E16").Value = 3.256523330529514
$ws.Range("H16").Value = 15.99231275138862
$ws.Range("I16").Value = 0.7883777257776909
$ws.Range("D17").Value = 18.83640741177098
$ws.Range("E17").Value = 2.027115801488069
$ws.Range("H17").Value = 17.26619296836839
$ws.Range("I17").Value = 0.4569013580854804
$ws.Range("D18").Value = 21.76884348995441
$ws.Range("E18").Value = 3.452233419874557
$ws.Range("H18").Value = 18.36869430832971
$ws.Range("I18").Value = 0.05208423824985831
$ws.Range("D19").Value = 16.84457784873161
$ws.Range("E19").Value = 0.0907765867928636
$ws.Range("H19").Value = 17.09805467329123
$ws.Range("I19").Value = 0.3442534113524888
$ws.Range("D20").Value = 16.73062730649282
$ws.Range("E20").Value = 1.077882637318507
$ws.Range("H20").Value = 16.10957976720388
$ws.Range("I20").Value = 0.4568350980295606
$ws.Range("D21").Value = 21.75321918691505
$ws.Range("E21").Value = 3.068422673434702
$ws.Range("H21").Value = 19.02011810422116
$ws.Range("I21").Value = 0.3353215907408094
$ws.Range("D22").Value = 19.74374860627168
$ws.Range("E22").Value = 2.26178141676997
$ws.Range("H22").Value = 17.79421490741208
$ws.Range("I22").Value = 0.312247717910368
$ws.Range("D23").Value = 21.55298493546261
$ws.Range("E23").Value = 4.222963177703192
$ws.Range("H23").Value = 17.82007041072784
$ws.Range("I23").Value = 0.4900486529684207
$ws.Range("D24").Value = 17.69860305879524
$ws.Range("E24").Value = 1.799156617467514
$ws.Range("H24").Value = 16.16013181305845
$ws.Range("I24").Value = 0.26068537173072
$ws.Range("D25").Value = 23.93557158630663
$ws.Range("E25").Value = 5.029055050542582
$ws.Range("H25").Value = 19.50912115187695
$ws.Range("I25").Value = 0.6026046161129013
$ws.Range("D26").Value = 13.1926505795087
$ws.Range("E26").Value = -2.250099503520514
$ws.Range("H26").Value = 15.49976509713668
$ws.Range("I26").Value = 0.05701501410746013
$ws.Range("D27").Value = 20.45247125591332
$ws.Range("E27").Value = 2.088654538747423
$ws.Range("H27").Value = 18.58144749615976
$ws.Range("I27").Value = 0.217630778993863
$ws.Range("D28").Value = 9.444197296833567
$ws.Range("E28").Value = -4.224344989696759
$ws.Range("H28").Value = 12.24625818749493
$ws.Range("I28").Value = -1.422284099035399
$ws.Range("D29").Value = 18.50946756792049
$ws.Range("E29").Value = 2.664548509501821
$ws.Range("H29").Value = 16.11308818455639
$ws.Range("I29").Value = 0.2681691261377193
$ws.Range("D30").Value = 20.14911725371923
$ws.Range("E30").Value = 2.513815827699128
$ws.Range("H30").Value = 17.88479765161451
$ws.Range("I30").Value = 0.2494962255944115
$ws.Range("D31").Value = 16.47611816641983
$ws.Range("E31").Value = 0.8630997574262131
$ws.Range("H31").Value = 15.68171524598895
$ws.Range("I31").Value = 0.06869683699533402
$ws.Range("D32").Value = 4.163539394345281
$ws.Range("E32").Value = -9.136259922448684
$ws.Range("H32").Value = 12.94094694019077
$ws.Range("I32").Value = -0.3588523766031944
$ws.Range("D33").Value = 6.700651309341816
$ws.Range("E33").Value = -5.965871328215023
$ws.Range("H33").Value = 13.09151700379999
$ws.Range("I33").Value = 0.4249943662431541
$ws.Range("D34").Value = 15.46849954260049
$ws.Range("E34").Value = -1.904601005127268
$ws.Range("H34").Value = 16.5400321226507
$ws.Range("I34").Value = -0.8330684250770529
$ws.Range("D35").Value = 1.290029740081234
$ws.Range("E35").Value = -10.2155807850813
$ws.Range("H35").Value = 10.23941986922571
$ws.Range("I35").Value = -1.266190655936818
$ws.Range("D36").Value = 8.869424641777156
$ws.Range("E36").Value = -5.118862888969086
$ws.Range("H36").Value = 14.38678012115603
$ws.Range("I36").Value = 0.3984925904097918
$ws.Range("D37").Value = -0.2179045357821263
$ws.Range("E37").Value = -12.03998891047272
$ws.Range("H37").Value = 11.66973432145659
$ws.Range("I37").Value = -0.1523500532340007
$ws.Range("D38").Value = 20.40555419301792
$ws.Range("E38").Value = 3.737216478130606
$ws.Range("H38").Value = 17.00934257144948
$ws.Range("I38").Value = 0.3410048565621724
$ws.Range("D39").Value = 5.457098949465291
$ws.Range("E39").Value = -7.610506153323709
$ws.Range("H39").Value = 13.00643176527852
$ws.Range("I39").Value = -0.06117333751047993
$ws.Range("D40").Value = 13.50545764681694
$ws.Range("E40").Value = -2.282902890863729
$ws.Range("H40").Value = 15.43035728033293
$ws.Range("I40").Value = -0.3580032573477372
$ws.Range("D41").Value = 0.4326074580196289
$ws.Range("E41").Value = -11.89321567174351
$ws.Range("H41").Value = 11.7810266597183
$ws.Range("I41").Value = -0.5447964700448378
$ws.Range("D42").Value = 2.189784869515821
$ws.Range("E42").Value = -11.40241170612937
$ws.Range("H42").Value = 11.62918995108879
$ws.Range("I42").Value = -1.963006624556405
$ws.Range("D43").Value = 9.444840389407188
$ws.Range("E43").Value = -4.825273439396264
$ws.Range("H43").Value = 14.35836068480814
$ws.Range("I43").Value = 0.08824685600468973
$ws.Range("D44").Value = 10.8966410157768
$ws.Range("E44").Value = -3.922149984252447
$ws.Range("H44").Value = 14.20980002551697
$ws.Range("I44").Value = -0.6089909745122815
$ws.Range("D45").Value = 4.677656171601808
$ws.Range("E45").Value = -9.269366245407104
$ws.Range("H45").Value = 13.5818674602828
$ws.Range("I45").Value = -0.365154956726105
$ws.Range("D46").Value = 9.632950178470395
$ws.Range("E46").Value = -3.252503968400178
$ws.Range("H46").Value = 12.77946002393264
$ws.Range("I46").Value = -0.1059941229379301
$ws.Range("D47").Value = -0.04312166362878056
$ws.Range("E47").Value = -11.8463498264206
$ws.Range("H47").Value = 10.35138281364132
$ws.Range("I47").Value = -1.451845349150501
$ws.Range("D48").Value = 5.450096855694061
$ws.Range("E48").Value = -7.39300016775293
$ws.Range("H48").Value = 12.5325814188397
$ws.Range("I48").Value = -0.3105156046072892
$ws.Range("D49").Value = 16.5712511836367
$ws.Range("E49").Value = 0.9687446217745062
$ws.Range("H49").Value = 15.62879530148034
$ws.Range("I49").Value = 0.02628873961814593
$ws.Range("D50").Value = 9.744475174887114
$ws.Range("E50").Value = -3.364140110546203
$ws.Range("H50").Value = 12.940556164606
$ws.Range("I50").Value = -0.168059120827315
$ws.Range("D51").Value = 3.19015705912469
$ws.Range("E51").Value = -9.988936645456667
$ws.Range("H51").Value = 12.65066150925681
$ws.Range("I51").Value = -0.5284321953245431
$ws.Range("D52").Value = 3.887200865728969
$ws.Range("E52").Value = -8.828320885929458
$ws.Range("H52").Value = 12.27047449964756
$ws.Range("I52").Value = -0.4450472520108657
$ws.Range("D53").Value = 8.068973851125318
$ws.Range("E53").Value = -6.021492597517403
$ws.Range("H53").Value = 12.96500470795793
$ws.Range("I53").Value = -1.125461740684786
$ws.Range("D54").Value = 13.81236353833143
$ws.Range("E54").Value = -0.5299936200610653
$ws.Range("H54").Value = 14.04025181893629
$ws.Range("I54").Value = -0.3021053394562006
$ws.Range("D55").Value = 7.721173906627175
$ws.Range("E55").Value = -6.148882735854456
$ws.Range("H55").Value = 13.37063521969007
$ws.Range("I55").Value = -0.4994214227915621
$ws.Range("D56").Value = 8.010112500213243
$ws.Range("E56").Value = -7.029057772125838
$ws.Range("H56").Value = 14.4656148288049
$ws.Range("I56").Value = -0.5735554435341843
$ws.Range("D57").Value = 14.27102271082316
$ws.Range("E57").Value = -1.857700269309519
$ws.Range("H57").Value = 15.63930846252772
$ws.Range("I57").Value = -0.4894145176049631
$ws.Range("D58").Value = 10.61860871541032
$ws.Range("E58").Value = -3.961569885558605
$ws.Range("H58").Value = 13.80617063219902
$ws.Range("I58").Value = -0.7740079687698995
$ws.Range("D59").Value = 12.38325472146277
$ws.Range("E59").Value = -3.593546478328274
$ws.Range("H59").Value = 15.51856988518849
$ws.Range("I59").Value = -0.4582313146025552
$ws.Range("D60").Value = 6.915033594442186
$ws.Range("E60").Value = -6.965133233092807
$ws.Range("H60").Value = 13.57406160196382
$ws.Range("I60").Value = -0.3061052255711711
$ws.Range("D61").Value = 10.065667974418
$ws.Range("E61").Value = -5.631213007999292
$ws.Range("H61").Value = 15.37282343184918
$ws.Range("I61").Value = -0.3240575505681202
$ws.Range("D62").Value = 3.883861942378701
$ws.Range("E62").Value = -8.008506176926371
$ws.Range("H62").Value = 10.36375654500825
$ws.Range("I62").Value = -1.528611574296823
$ws.Range("D63").Value = 7.459984682547423
$ws.Range("E63").Value = -6.074528113362954
$ws.Range("H63").Value = 13.21218344617723
$ws.Range("I63").Value = -0.3223293497331468
$ws.Range("D64").Value = 8.529283775704425
$ws.Range("E64").Value = -4.058390986055885
$ws.Range("H64").Value = 12.16555603512687
$ws.Range("I64").Value = -0.4221187266334434
$ws.Range("D65").Value = 10.26670366666087
$ws.Range("E65").Value = -3.182457534955237
$ws.Range("H65").Value = 12.81889490196887
$ws.Range("I65").Value = -0.6302662996472309
$ws.Range("D66").Value = 12.14539756530341
$ws.Range("E66").Value = -3.500541476163075
$ws.Range("H66").Value = 14.74587131904596
$ws.Range("I66").Value = -0.9000677224205231
$ws.Range("D67").Value = 14.25548349653589
$ws.Range("E67").Value = -2.144863523659787
$ws.Range("H67").Value = 16.03950887060596
$ws.Range("I67").Value = -0.3608381495897197
$ws.Range("D68").Value = 10.65278613007859
$ws.Range("E68").Value = -4.540562680606897
$ws.Range("H68").Value = 15.00641625846313
$ws.Range("I68").Value = -0.1869325522223557
$ws.Range("D69").Value = 0.7282636210301494
$ws.Range("E69").Value = -12.37847365879339
$ws.Range("H69").Value = 10.92716328375705
$ws.Range("I69").Value = -2.17957399606648
$ws.Range("D70").Value = 2.282635128447745
$ws.Range("E70").Value = -10.26859236748521
$ws.Range("H70").Value = 11.5925298019141
$ws.Range("I70").Value = -0.9586976940188537
$ws.Range("D71").Value = 5.330430744524403
$ws.Range("E71").Value = -6.833441600343041
$ws.Range("H71").Value = 12.15589795031265
$ws.Range("I71").Value = -0.007974394554788544
$ws.Range("D72").Value = 7.244830263202964
$ws.Range("E72").Value = -6.096046700624055
$ws.Range("H72").Value = 13.46687605617228
$ws.Range("I72").Value = 0.1259990923452587
$ws.Range("D73").Value = 3.937704483914213
$ws.Range("E73").Value = -10.92096978930029
$ws.Range("H73").Value = 13.70576338357865
$ws.Range("I73").Value = -1.152910889635853
$ws.Range("D74").Value = 15.56620351631875
$ws.Range("E74").Value = -0.2747738524513517
$ws.Range("H74").Value = 15.80579192396723
$ws.Range("I74").Value = -0.03518544480287389
$ws.Range("D75").Value = 9.555375105847418
$ws.Range("E75").Value = -5.635557377856915
$ws.Range("H75").Value = 14.14360699219884
$ws.Range("I75").Value = -1.047325491505492
$ws.Range("D76").Value = 4.235131808529415
$ws.Range("E76").Value = -9.369771775333938
$ws.Range("H76").Value = 13.02935014981891
$ws.Range("I76").Value = -0.5755534340444481
$ws.Range("D77").Value = 21.81993832794947
$ws.Range("E77").Value = 4.850394062388732
$ws.Range("H77").Value = 17.38532177877094
$ws.Range("I77").Value = 0.4157775132102051
$ws.Range("D78").Value = 6.259807198565517
$ws.Range("E78").Value = -6.951768133970831
$ws.Range("H78").Value = 13.436752347535
$ws.Range("I78").Value = 0.2251770149986563
$ws.Range("D79").Value = -4.562035244427765
$ws.Range("E79").Value = -15.18220275083886
$ws.Range("H79").Value = 10.50316770751554
$ws.Range("I79").Value = -0.1169997988955604
$ws.Range("D80").Value = 10.53257481703519
$ws.Range("E80").Value = -4.864926082841744
$ws.Range("H80").Value = 14.67246039086362
$ws.Range("I80").Value = -0.7250405090133138
$ws.Range("D81").Value = 10.08473596218145
$ws.Range("E81").Value = -3.882452654211812
$ws.Range("H81").Value = 12.92074549033922
$ws.Range("I81").Value = -1.046443126054043
$ws.Range("D82").Value = 3.315958991411651
$ws.Range("E82").Value = -8.899602710754374
$ws.Range("H82").Value = 11.5956129556812
$ws.Range("I82").Value = -0.6199487464848215
$ws.Range("D83").Value = 0.1259356123694182
$ws.Range("E83").Value = -10.65962315711549
$ws.Range("H83").Value = 10.84303168236413
$ws.Range("I83").Value = 0.05747291287922274
$ws.Range("D84").Value = 19.13239754639152
$ws.Range("E84").Value = 1.623053620690789
$ws.Range("H84").Value = 17.52600308348096
$ws.Range("I84").Value = 0.01665915778023397
$ws.Range("D85").Value = 20.64329527154715
$ws.Range("E85").Value = 4.846593031871143
$ws.Range("H85").Value = 15.99686090380953
$ws.Range("I85").Value = 0.2001586641335227
$ws.Range("D86").Value = 19.44321956252464
$ws.Range("E86").Value = 2.679982560502483
$ws.Range("H86").Value = 17.51814402493251
$ws.Range("I86").Value = 0.7549070229103636
$ws.Range("D87").Value = 19.21146202954246
$ws.Range("E87").Value = 3.646604166417159
$ws.Range("H87").Value = 16.03268281645457
$ws.Range("I87").Value = 0.4678249533292647
$ws.Range("D88").Value = 20.78428220439769
$ws.Range("E88").Value = 4.5610573049581
$ws.Range("H88").Value = 16.72153745380217
$ws.Range("I88").Value = 0.4983125543625782
$ws.Range("D89").Value = 22.55241826879008
$ws.Range("E89").Value = 4.028996707839013
$ws.Range("H89").Value = 18.82696948754088
$ws.Range("I89").Value = 0.3035479265898129
$ws.Range("D90").Value = 23.54790193612411
$ws.Range("E90").Value = 4.587776230702331
$ws.Range("H90").Value = 19.39370739571533
$ws.Range("I90").Value = 0.4335816902935463
$ws.Range("D91").Value = 20.46695203319463
$ws.Range("E91").Value = 4.677282950777613
$ws.Range("H91").Value = 16.15454494125146
$ws.Range("I91").Value = 0.3648758588344401
$ws.Range("D92").Value = 19.31299628409622
$ws.Range("E92").Value = 3.551075877309152
$ws.Range("H92").Value = 16.23151027051814
$ws.Range("I92").Value = 0.4695898637310767
$ws.Range("D93").Value = 19.3758177425279
$ws.Range("E93").Value = 2.321334824279698
$ws.Range("H93").Value = 17.38233902370263
$ws.Range("I93").Value = 0.3278561054544376
$ws.Range("D94").Value = 19.84008398578937
$ws.Range("E94").Value = 3.483511751659217
$ws.Range("H94").Value = 16.36430570508444
$ws.Range("I94").Value = 0.007733470954281341
$ws.Range("D95").Value = 19.13388324359203
$ws.Range("E95").Value = 4.169268552535679
$ws.Range("H95").Value = 15.4076824667096
$ws.Range("I95").Value = 0.4430677756532504
$ws.Range("D96").Value = 16.72898384748754
$ws.Range("E96").Value = 0.09213226907106631
$ws.Range("H96").Value = 16.99333791382486
$ws.Range("I96").Value = 0.3564863354083869
$ws.Range("D97").Value = 19.22345797784778
$ws.Range("E97").Value = 2.565876245399097
$ws.Range("H97").Value = 17.1229906934518
$ws.Range("I97").Value = 0.4654089610031187
$ws.Range("D98").Value = 21.2118843362794
$ws.Range("E98").Value = 3.171008216561439
$ws.Range("H98").Value = 18.23536628934741
$ws.Range("I98").Value = 0.194490169629439
$ws.Range("D99").Value = 20.67691923524711
$ws.Range("E99").Value = 3.171549220472687
$ws.Range("H99").Value = 17.95338632787417
$ws.Range("I99").Value = 0.4480163130997457
$ws.Range("D100").Value = 21.39912884983688
$ws.Range("E100").Value = 2.595683755447869
$ws.Range("H100").Value = 19.04654748003089
$ws.Range("I100").Value = 0.2431023856418781
$ws.Range("D101").Value = 21.88605803047528
$ws.Range("E101").Value = 2.935972141374082
$ws.Range("H101").Value = 19.09102926709085
$ws.Range("I101").Value = 0.1409433779896487
$ws.Range("D102").Value = 20.6541965467521
$ws.Range("E102").Value = 4.101958489664043
$ws.Range("H102").Value = 16.96572613193494
$ws.Range("I102").Value = 0.4134880748468792
$ws.Range("D103").Value = 17.00594582607346
$ws.Range("E103").Value = 1.938154385302206
$ws.Range("H103").Value = 15.23897831426211
$ws.Range("I103").Value = 0.1711868734908519
$ws.Range("D104").Value = 23.07494644811731
$ws.Range("E104").Value = 5.146573681899988
$ws.Range("H104").Value = 18.64527077863276
$ws.Range("I104").Value = 0.716898012415444
$ws.Range("D105").Value = 22.30755258051497
$ws.Range("E105").Value = 4.485906415790654
$ws.Range("H105").Value = 18.03582665341589
$ws.Range("I105").Value = 0.2141804886915766
$ws.Range("D106").Value = 23.33035276225481
$ws.Range("E106").Value = 5.173219317508185
$ws.Range("H106").Value = 18.47081797856139
$ws.Range("I106").Value = 0.3136845338147669
$ws.Range("D107").Value = 19.62428675504104
$ws.Range("E107").Value = 3.191826912370783
$ws.Range("H107").Value = 16.5889343981664
$ws.Range("I107").Value = 0.1564745554961409
$ws.Range("D108").Value = 20.18999005460871
$ws.Range("E108").Value = 2.831260834546474
$ws.Range("H108").Value = 17.90890454081421
$ws.Range("I108").Value = 0.5501753207519751
$ws.Range("D109").Value = 23.30180004506395
$ws.Range("E109").Value = 4.30621001710618
$ws.Range("H109").Value = 19.50662786232684
$ws.Range("I109").Value = 0.5110378343690656
$ws.Range("D110").Value = 20.81188561803199
$ws.Range("E110").Value = 4.63277368018231
$ws.Range("H110").Value = 16.62937956470086
$ws.Range("I110").Value = 0.4502676268511783
$ws.Range("D112").Value = 24.46875725844992
$ws.Range("E112").Value = 5.588197887207379
$ws.Range("H112").Value = 19.59677427516738
$ws.Range("I112").Value = 0.7162149039248378
$ws.Range("D113").Value = 21.8406189407292
$ws.Range("E113").Value = 3.607802218206481
$ws.Range("H113").Value = 18.93052579815172
$ws.Range("I113").Value = 0.6977090756290052
$ws.Range("D114").Value = 23.11455553150648
$ws.Range("E114").Value = 4.88153447821975
$ws.Range("H114").Value = 18.69544667164335
$ws.Range("I114").Value = 0.4624256183566264
$ws.Range("D115").Value = 25.36631015521095
$ws.Range("E115").Value = 6.565137294807529
$ws.Range("H115").Value = 19.35185659036009
$ws.Range("I115").Value = 0.5506837299566618
$ws.Range("D116").Value = 21.15179991290754
$ws.Range("E116").Value = 3.01808161245088
$ws.Range("H116").Value = 18.95188237548249
$ws.Range("I116").Value = 0.8181640750258272
$ws.Range("D117").Value = 19.74574237478612
$ws.Range("E117").Value = 3.869653624270993
$ws.Range("H117").Value = 16.24177858861682
$ws.Range("I117").Value = 0.3656898381016986
$ws.Range("D118").Value = 24.93296375059332
$ws.Range("E118").Value = 6.858895542324948
$ws.Range("H118").Value = 18.65359586628811
$ws.Range("I118").Value = 0.5795276580197419
$ws.Range("D119").Value = 25.1202082641508
$ws.Range("E119").Value = 6.283571081211377
$ws.Range("H119").Value = 19.4647770569716
$ws.Range("I119").Value = 0.6281398740321811
$ws.Range("D120").Value = 20.32075356649137
$ws.Range("E120").Value = 4.62196378127596
$ws.Range("H120").Value = 15.9623159882125
$ws.Range("I120").Value = 0.2635262029970886
$ws.Range("D121").Value = 24.01210842826751
$ws.Range("E121").Value = 5.8584738858199
$ws.Range("H121").Value = 18.45052898683606
$ws.Range("I121").Value = 0.2968944443884505
$ws.Range("D122").Value = 23.53086624973617
$ws.Range("E122").Value = 4.104717983284144
$ws.Range("H122").Value = 19.97212414208738
$ws.Range("I122").Value = 0.5459758756353488
$ws.Range("D123").Value = 18.40739301976334
$ws.Range("E123").Value = 2.967669785126604
$ws.Range("H123").Value = 15.85862401679311
$ws.Range("I123").Value = 0.4189007821563628
$ws.Range("D124").Value = 23.57876202364988
$ws.Range("E124").Value = 6.15223213333732
$ws.Range("H124").Value = 17.75226826276409
$ws.Range("I124").Value = 0.3257383724515306
$ws.Range("D125").Value = 23.90661999573352
$ws.Range("E125").Value = 4.735660116826642
$ws.Range("H125").Value = 19.7282970910649
$ws.Range("I125").Value = 0.5573372121580242
$ws.Range("D126").Value = 15.52304868908039
$ws.Range("E126").Value = 0.4208253198561652
$ws.Range("H126").Value = 15.11004505196338
$ws.Range("I126").Value = 0.007821682739157954
$ws.Range("D127").Value = 22.63710681951976
$ws.Range("E127").Value = 5.625332732610541
$ws.Range("H127").Value = 18.21188787173196
$ws.Range("I127").Value = 1.200113784822742
